# Apply updated crypto price/volume figures to Sheet1 (cryptos.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.426.89'
$ws.Range('E2').Value = '  +0.39%  '
$ws.Range('D3').Value = '2.931.64'
$ws.Range('E3').Value = '  -0.03%  '
$ws.Range('E4').Value = '  +0.03%  '
$r = $ws.Range('D5')
$r.NumberFormat = "@"
$r.Value = '598.23'
$r.Style = "Normal"
$ws.Range('E5').Value = '  +1.05%  '
$r = $ws.Range('D6')
$r.NumberFormat = "@"
$r.Value = '145.51'
$r.Style = "Normal"
$ws.Range('E6').Value = '  -0.70%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -0.90%  '
$r = $ws.Range('D9')
$r.NumberFormat = "@"
$r.Value = '7.02'
$r.Style = "Normal"
$ws.Range('E9').Value = '  +1.55%  '
$ws.Range('E10').Value = '  -2.01%  '
$ws.Range('E11').Value = '  -0.72%  '
$ws.Range('E12').Value = '  -1.15%  '
$r = $ws.Range('D13')
$r.NumberFormat = "@"
$r.Value = '33.56'
$r.Style = "Normal"
$ws.Range('E13').Value = '  -0.91%  '
$ws.Range('E14').Value = '  +0.39%  '
$ws.Range('D15').Value = '3.418.35'
$ws.Range('E15').Value = '  +0.04%  '
$ws.Range('D16').Value = '61.413.78'
$ws.Range('E16').Value = '  +0.49%  '
$ws.Range('D17').Value = '2.931.92'
$ws.Range('E17').Value = '  +0.47%  '
$ws.Range('E18').Value = '  -0.39%  '
$r = $ws.Range('D19')
$r.NumberFormat = "@"
$r.Value = '431.90'
$r.Style = "Normal"
$ws.Range('E19').Value = '  -0.11%  '
$r = $ws.Range('D20')
$r.NumberFormat = "@"
$r.Value = '13.47'
$r.Style = "Normal"
$ws.Range('E20').Value = '  +0.02%  '
$ws.Range('E21').Value = '  -1.25%  '
$r = $ws.Range('D22')
$r.NumberFormat = "@"
$r.Value = '7.09'
$r.Style = "Normal"
$ws.Range('E22').Value = '  -0.27%  '
$r = $ws.Range('D23')
$r.NumberFormat = "@"
$r.Value = '81.96'
$r.Style = "Normal"
$ws.Range('E23').Value = '  +0.67%  '
$r = $ws.Range('D24')
$r.NumberFormat = "@"
$r.Value = '10.88'
$r.Style = "Normal"
$ws.Range('E24').Value = '  -1.76%  '
$ws.Range('E25').Value = '  -2.15%  '
$ws.Range('E26').Value = '  -2.14%  '
$ws.Range('E27').Value = '  -0.01%  '
$r = $ws.Range('D28')
$r.NumberFormat = "@"
$r.Value = '2.22'
$r.Style = "Normal"
$ws.Range('E28').Value = '  -4.63%  '
$r = $ws.Range('D29')
$r.NumberFormat = "@"
$r.Value = '2.61'
$r.Style = "Normal"
$ws.Range('E29').Value = '  -0.33%  '
$r = $ws.Range('D30')
$r.NumberFormat = "@"
$r.Value = '6.93'
$r.Style = "Normal"
$ws.Range('E30').Value = '  -2.71%  '
$r = $ws.Range('D31')
$r.NumberFormat = "@"
$r.Value = '26.63'
$r.Style = "Normal"
$ws.Range('E31').Value = '  +0.36%  '
$ws.Range('E32').Value = '  +1.00%  '
$ws.Range('E33').Value = '  +0.13%  '
$ws.Range('D34').Value = '0.0₃0884'
$ws.Range('E34').Value = '  +2.29%  '
$ws.Range('E35').Value = '  -0.23%  '
$ws.Range('E36').Value = '  -0.09%  '
$r = $ws.Range('D37')
$r.NumberFormat = "@"
$r.Value = '2.98'
$r.Style = "Normal"
$ws.Range('E37').Value = '  -3.57%  '
$r = $ws.Range('D38')
$r.NumberFormat = "@"
$r.Value = '2.00'
$r.Style = "Normal"
$ws.Range('E38').Value = '  -0.26%  '
$ws.Range('E39').Value = '  -1.73%  '
$ws.Range('E40').Value = '  -0.30%  '
$ws.Range('E41').Value = '  +6.49%  '
$ws.Range('E42').Value = '  -2.50%  '
$r = $ws.Range('D43')
$r.NumberFormat = "@"
$r.Value = '0.0346'
$r.Style = "Normal"
$ws.Range('E43').Value = '  -0.30%  '
$ws.Range('D44').Value = '2.697.15'
$ws.Range('E44').Value = '  -0.73%  '
$r = $ws.Range('D45')
$r.NumberFormat = "@"
$r.Value = '134.23'
$r.Style = "Normal"
$ws.Range('E45').Value = '  +2.06%  '
$r = $ws.Range('D46')
$r.NumberFormat = "@"
$r.Value = '361.03'
$r.Style = "Normal"
$ws.Range('E46').Value = '  -4.03%  '
$r = $ws.Range('D48')
$r.NumberFormat = "@"
$r.Value = '23.65'
$r.Style = "Normal"
$ws.Range('E48').Value = '  -2.74%  '
$ws.Range('E49').Value = '  -1.36%  '
$ws.Range('E50').Value = '  -2.03%  '
$ws.Range('E51').Value = '  -1.61%  '
